$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Form Responses 1")

# --- Row 2: update existing respondent's details ---
$ws.Range("B2").Value = "Himanshu Kashyap"
$ws.Range("C2").Value = "MU19BTCSE017"
$ws.Range("D2").Value = "himmu7987@gmail.com"
$ws.Range("E2").Value = "F"
$ws.Range("E2").NumberFormat = "General"

# --- Row 3: update existing respondent's details ---
$ws.Range("B3").Value = "Prakash Singh"
$ws.Range("C3").Value = "MU19BTCSE002L"
$ws.Range("D3").Value = "kshatreeya545@gmail.com"
$ws.Range("E3").Value = "T"
$ws.Range("E3").NumberFormat = "General"
$ws.Range("A3:E3").Font.Name = "Calibri"

# --- Header row cleanup: Certificate Status column format ---
$ws.Range("E1").NumberFormat = "General"

# --- Row 4: remove the third (now-superfluous) submission ---
$ws.Range("A4:E4").Delete()
